# Add "Deferred gross" and "Taxable" columns to the Detail sheet of the
# Account Payroll template, immediately after the existing "Hours Worked"
# column (K = Def Gross, L = Taxable), matching the existing bold/centered
# header formatting used by the rest of row 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Detail")

# New header text - lands in xl/sharedStrings.xml and row 1 of the sheet.
$ws.Range("K1").Value = "Def Gross"
$ws.Range("L1").Value = "Taxable"

# Match the header formatting (bold "Heading 2" style with bottom border)
# used by the other header cells, by copying the format from the adjacent
# "Hours Worked" header cell instead of re-creating a new style.
$ws.Range("J1").Copy()
$ws.Range("K1:L1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
